# Update models.xlsx to reflect models used to make predictions.
#
# The "Predicted?" column (I) on the "all" sheet is updated for a handful
# of rows: some models that were previously marked "No"/"No "/"in progress"
# are now marked "Yes" (predictions made) or "In progress".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("all")

# Row 2: SEStructure / bin_GFP -> prediction finished
$ws.Cells.Item(2, 9).Value = "Yes"

# Row 3: 32Structure 52Structure / bin_GFP -> prediction finished
$ws.Cells.Item(3, 9).Value = "Yes"

# Row 4: 52Structure / bin_GFP -> prediction finished
$ws.Cells.Item(4, 9).Value = "Yes"

# Row 11: structure / bin_mKate -> prediction finished
$ws.Cells.Item(11, 9).Value = "Yes"

# Row 12: SEStructure / bin_mKate -> prediction in progress
$ws.Cells.Item(12, 9).Value = "In progress"

# Row 13: 32Structure 52Structure / bin_mKate -> prediction in progress
$ws.Cells.Item(13, 9).Value = "In progress"

# Row 27: 52Structure / log_gfp -> prediction in progress
$ws.Cells.Item(27, 9).Value = "In progress"

# Row 33: structure / log_mKate -> prediction in progress
$ws.Cells.Item(33, 9).Value = "In progress"

# Row 36: 52Structure / log_mKate -> prediction in progress
$ws.Cells.Item(36, 9).Value = "In progress"

# Update the view: scrolled down to show row 12 at top, with I36 selected
# (reflecting the row the author was last looking at/editing).
$ws.Range("I36").Select()
try {
    $excel.ActiveWindow.ScrollRow = 12
    $excel.ActiveWindow.ScrollColumn = 1
} catch {
    # Scroll position is a cosmetic view setting; ignore if unsupported.
}
